$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 24; $r -le 248; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $r - 2
}
